# Update "想去人数" (F column) figures on the sheets that list individual
# events ("展览" and "全部类型"). Both sheets share the same event rows, so
# the same updates are applied to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1379
    3 = 2183
    4 = 325
    5 = 76
    6 = 6409
    7 = 282
    8 = 120
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
